$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L14").Value = -33.333333333333

$ws.Range("N14").Value = -42.857142857142

$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("L15").Value = -16.666666666666

$ws.Range("N15").Value = -37.5

$ws.Range("C16").Value = 2

$ws.Range("D16").Value = 3

$ws.Range("E16").Value = -33.333333333333

$ws.Range("F16").Value = 10

$ws.Range("H16").Value = 11.111111111111

$ws.Range("I16").Value = 78

$ws.Range("J16").Value = 72

$ws.Range("K16").Value = 8.333333333333

$ws.Range("L16").Value = 200

$ws.Range("M16").Value = -38.582677165354

$ws.Range("N16").Value = -78.688524590163

$ws.Range("C17").Value = 8

$ws.Range("D17").Value = 3

$ws.Range("E17").Value = 166.666666666667

$ws.Range("F17").Value = 21

$ws.Range("G17").Value = 14

$ws.Range("H17").Value = 50

$ws.Range("I17").Value = 137

$ws.Range("J17").Value = 141

$ws.Range("K17").Value = -2.836879432624

$ws.Range("L17").Value = 28.03738317757

$ws.Range("M17").Value = 44.210526315789

$ws.Range("N17").Value = -11.612903225806

$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C18").PasteSpecial(-4122)

$ws.Range("E18").Value = -100

$ws.Range("F18").Value = 3

$ws.Range("H18").Value = -66.666666666666

$ws.Range("J18").Value = 60

$ws.Range("K18").Value = -46.666666666666

$ws.Range("M18").Value = -78.666666666666

$ws.Range("N18").Value = -91.280653950953

$ws.Range("D19").Value = 6

$ws.Range("E19").Value = -16.666666666666

$ws.Range("F19").Value = 15

$ws.Range("G19").Value = 20

$ws.Range("H19").Value = -25

$ws.Range("I19").Value = 127

$ws.Range("J19").Value = 152

$ws.Range("K19").Value = -16.447368421052

$ws.Range("L19").Value = 58.75

$ws.Range("M19").Value = -2.307692307692

$ws.Range("N19").Value = -46.186440677966

$ws.Range("C20").Value = 3

$ws.Range("D20").Value = 5

$ws.Range("E20").Value = -40

$ws.Range("F20").Value = 12

$ws.Range("G20").Value = 12

$ws.Range("H20").Value = 0

$ws.Range("I20").Value = 94

$ws.Range("J20").Value = 89

$ws.Range("K20").Value = 5.617977528089

$ws.Range("L20").Value = 74.074074074074

$ws.Range("M20").Value = 40.298507462686

$ws.Range("N20").Value = -92.503987240829

$ws.Range("C21").Value = 18

$ws.Range("D21").Value = 19

$ws.Range("E21").Value = -5.263157894736

$ws.Range("F21").Value = 64

$ws.Range("G21").Value = 66

$ws.Range("H21").Value = -3.030303030303

$ws.Range("I21").Value = 482

$ws.Range("J21").Value = 530

$ws.Range("K21").Value = -9.056603773584

$ws.Range("L21").Value = 47.852760736196

$ws.Range("M21").Value = -17.324185248713

$ws.Range("N21").Value = -79.925031236984

$ws.Range("C23").Value = 2

$ws.Range("D23").Value = 5

$ws.Range("E23").Value = -60

$ws.Range("G23").Value = 13

$ws.Range("H23").Value = -23.076923076923

$ws.Range("I23").Value = 75

$ws.Range("J23").Value = 64

$ws.Range("K23").Value = 17.1875

$ws.Range("L23").Value = 33.928571428571

$ws.Range("M23").Value = 141.935483870968

$ws.Range("C24").Value = 13

$ws.Range("E24").Value = -23.529411764705

$ws.Range("F24").Value = 49

$ws.Range("G24").Value = 50

$ws.Range("H24").Value = -2

$ws.Range("I24").Value = 354

$ws.Range("J24").Value = 328

$ws.Range("K24").Value = 7.926829268292

$ws.Range("L24").Value = 46.280991735537

$ws.Range("M24").Value = 7.926829268292

$ws.Range("C25").Value = 7

$ws.Range("D25").Value = 7

$ws.Range("E25").Value = 0

$ws.Range("F25").Value = 32

$ws.Range("G25").Value = 31

$ws.Range("H25").Value = 3.225806451612

$ws.Range("I25").Value = 228

$ws.Range("J25").Value = 235

$ws.Range("K25").Value = -2.978723404255

$ws.Range("L25").Value = 35.714285714285

$ws.Range("M25").Value = -22.711864406779

$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C26").PasteSpecial(-4122)

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D26").PasteSpecial(-4122)

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E26").PasteSpecial(-4122)

$ws.Range("L26").Value = -23.529411764705

$ws.Range("G27").Value = 1

$ws.Range("H27").Value = 200

$ws.Range("L27").Value = -17.647058823529

$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4122)

$ws.Range("G14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = 1

$ws.Range("H14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = -100

$ws.Range("G28").Value = 3

$ws.Range("H28").Value = 0

$ws.Range("J28").Value = 19

$ws.Range("K28").Value = -36.842105263157

$ws.Range("L28").Value = -36.842105263157

$ws.Range("M28").Value = -20

$ws.Range("N28").Value = -45.454545454545

$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C29").PasteSpecial(-4122)

$ws.Range("G14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D29").Value = 1

$ws.Range("H14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").Value = -100

$ws.Range("G29").Value = 3

$ws.Range("H29").Value = 0

$ws.Range("J29").Value = 16

$ws.Range("K29").Value = -25

$ws.Range("L29").Value = -20

$ws.Range("M29").Value = 0

$ws.Range("N29").Value = -25

$ws.Range("A8").Value = "Volume 30   Number  31"
$ws.Range("C9").Value = "Report Covering the Week  7/31/2023  Through  8/6/2023"